# Dictionnaire de donnees - ajout de la "longueur max" pour la table "facture"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("facture")

# Renseigne le champ "type" (D4) qui etait vide -> "string" (vehicule / FK)
$ws.Range("D4").Value = "string"

# Renseigne la colonne "longueur max" (C) avec un format numerique a separateur de milliers
$ws.Range("C2:C3").NumberFormat = "#,##0"
$ws.Range("C5:C6").NumberFormat = "#,##0"

$ws.Range("C2").Value = 6
$ws.Range("C3").Value = 14
$ws.Range("C5").Value = 6
$ws.Range("C6").Value = 6

# Retour sur l'onglet "vehicule" (premier onglet) comme onglet actif
$wb.Worksheets.Item("vehicule").Activate()
